$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently ends with:
#   row 37 -> totals row   (P37/Q37 merged, style 13)
#   row 38 -> footer row   (A38:F38 / G38:I38 / K38:Q38 merged, styles 14-17)
#
# The new version inserts one more sale-line (row 37, "31"), which pushes the
# totals row down to 38 (with an updated total) and the footer row down to 39.
#
# To keep the existing cell styles (7-17) exactly as they are (and avoid
# Excel fabricating brand-new style records), we never use Rows.Insert().
# Instead we grow into the previously untouched row 39, then shuffle the
# *formatting* downward one row at a time using copy / PasteSpecial
# (xlPasteFormats), and only afterwards overwrite the cell values. Because
# PasteSpecial only touches formatting, the values already sitting in the
# destination rows survive each formatting copy untouched.
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122

# 1) Copy the footer row's formatting (row 38) down into the brand-new row 39.
$ws.Range("A38:Q38").Copy()
$ws.Range("A39:Q39").PasteSpecial($xlPasteFormats)

# 2) Copy the (current) totals row's formatting (row 37) down into row 38.
$ws.Range("A37:Q37").Copy()
$ws.Range("A38:Q38").PasteSpecial($xlPasteFormats)

# 3) Copy a normal sale-line row's formatting (row 36) down into row 37, which
#    turns row 37 into a regular data row.
$ws.Range("A36:Q36").Copy()
$ws.Range("A37:Q37").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Undo the merges that used to sit on rows 37/38 before we rebuild them on
# their new rows.
# ---------------------------------------------------------------------------
$ws.Range("P37:Q37").UnMerge()
$ws.Range("A38:F38").UnMerge()
$ws.Range("G38:I38").UnMerge()
$ws.Range("K38:Q38").UnMerge()

# ---------------------------------------------------------------------------
# Row 39: footer (timestamp / page / developer credit), values updated.
# ---------------------------------------------------------------------------
$ws.Range("A39").Value = "Saturday, 2 August, 2025 1:31 PM"
$ws.Range("G39").Value = "1/1"
$ws.Range("K39").Value = "developed by : Abdelaziz Talaat"
$ws.Range("A39:F39").Merge()
$ws.Range("G39:I39").Merge()
$ws.Range("K39:Q39").Merge()
$ws.Rows.Item(39).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Row 38: totals row, with the new grand total.
# ---------------------------------------------------------------------------
$ws.Range("P38").Value = 1722.26
$ws.Range("Q38").Value = ""
$ws.Range("P38:Q38").Merge()
$ws.Rows.Item(38).RowHeight = 24.75

# ---------------------------------------------------------------------------
# Row 37: the new sale line (#31).
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = 31
$ws.Range("C37").Value = "كريم فيبكس الازرق"
$ws.Range("H37").Value = "2:0"
$ws.Range("L37").Value = "0"
$ws.Range("N37").Value = "35.00"
$ws.Range("P37").Value = "35.0000"
$ws.Range("Q37").Value = "1:0"

$ws.Range("A37:B37").Merge()
$ws.Range("C37:G37").Merge()
$ws.Range("H37:K37").Merge()
$ws.Range("L37:M37").Merge()
$ws.Range("N37:O37").Merge()
$ws.Rows.Item(37).RowHeight = 25.5

Write-Host "done"
